# Edit script: insert 2 new data rows (145 and 146) into Sheet1, pushing
# existing rows 145..216 down to 147..218, and populate the two new rows
# with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 145. This shifts the
# existing rows 145..216 down to 147..218 (matching dimension A1:R218).
$ws.Rows("145:146").Insert()

# --- New row 145 ---
$ws.Range("A145").Value = 3
$ws.Range("B145").Value = "Femacal de La Calera"
$ws.Range("C145").Value = "Coquimbo"
$ws.Range("D145").Value = 44489
$ws.Range("E145").Value = 5
$ws.Range("F145").Value = 100112013
$ws.Range("G145").Value = "Alcachofa"
$ws.Range("H145").Value = "Española"
$ws.Range("I145").Value = "Extra"
$ws.Range("J145").Value = 5800
$ws.Range("K145").Value = 380
$ws.Range("L145").Value = 380
$ws.Range("M145").Value = 380
$ws.Range("N145").Value = "$/unidad"
$ws.Range("O145").Value = "Llay Llay"
$ws.Range("P145").Value = 380
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = "Hortaliza"

# --- New row 146 ---
$ws.Range("A146").Value = 3
$ws.Range("B146").Value = "Femacal de La Calera"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44489
$ws.Range("E146").Value = 5
$ws.Range("F146").Value = 100112013
$ws.Range("G146").Value = "Alcachofa"
$ws.Range("H146").Value = "Española"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 5500
$ws.Range("K146").Value = 300
$ws.Range("L146").Value = 300
$ws.Range("M146").Value = 300
$ws.Range("N146").Value = "$/unidad"
$ws.Range("O146").Value = "Llay Llay"
$ws.Range("P146").Value = 300
$ws.Range("Q146").Value = 1
$ws.Range("R146").Value = "Hortaliza"
